$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "22.52000000000008"
$ws.Range("H2").Value = "1.709743457922741e-13"
$ws.Range("I2").Value = "1.709743457922741e-13"
$ws.Range("L2").Value = "44.67160608855833"
$ws.Range("M2").Value = "[35.39195167819982, 53.95126049891684]"
$ws.Range("N2").Value = "1.356470491487016e-12"
$ws.Range("O2").Value = "1.356470491487016e-12"
$ws.Range("P2").Value = "1.691868716347656"
$ws.Range("Q2").Value = "[1.478026573760964, 1.9057108589343486]"
$ws.Range("T2").Value = "52.22448494165594"
$ws.Range("U2").Value = "[46.581885202890795, 57.867084680421094]"
$ws.Range("X2").Value = "16.45605605605612"
$ws.Range("Y2").Value = "15.68960960960967"
$ws.Range("Z2").Value = "17.22250250250256"
$ws.Range("F3").Value = "22.52000000000008"
$ws.Range("H3").Value = "2.220446049250313e-16"
$ws.Range("I3").Value = "2.220446049250313e-16"
$ws.Range("L3").Value = "49.76888184634021"
$ws.Range("M3").Value = "[39.9677705291899, 59.569993163490516]"
$ws.Range("N3").Value = "2.557953848736361e-13"
$ws.Range("O3").Value = "2.557953848736361e-13"
$ws.Range("P3").Value = "2.056658253701427"
$ws.Range("Q3").Value = "[1.8428161111147352, 2.270500396288118]"
$ws.Range("T3").Value = "53.80613264221489"
$ws.Range("U3").Value = "[48.510999371608165, 59.10126591282162]"
$ws.Range("X3").Value = "15.14858858858864"
$ws.Range("Y3").Value = "14.3821421421422"
$ws.Range("Z3").Value = "15.91503503503509"
$ws.Range("F4").Value = "22.52000000000008"
$ws.Range("L4").Value = "48.06853571997853"
$ws.Range("M4").Value = "[38.28415282017755, 57.852918619779516]"
$ws.Range("N4").Value = "7.23421322845752e-13"
$ws.Range("O4").Value = "7.23421322845752e-13"
$ws.Range("P4").Value = "2.39628989192735"
$ws.Range("Q4").Value = "[2.19502669890458, 2.5975530849501203]"
$ws.Range("T4").Value = "46.49974887150642"
$ws.Range("U4").Value = "[41.54691336834092, 51.45258437467192]"
$ws.Range("X4").Value = "13.93129129129134"
$ws.Range("Y4").Value = "13.20992992992998"
$ws.Range("Z4").Value = "14.65265265265271"
$ws.Range("F5").Value = "23.11000000000017"
$ws.Range("H5").Value = "1.110223024625157e-16"
$ws.Range("I5").Value = "1.110223024625157e-16"
$ws.Range("L5").Value = "53.37011042452857"
$ws.Range("M5").Value = "[43.20379181881499, 63.53642903024215]"
$ws.Range("N5").Value = "8.815170815523743e-14"
$ws.Range("O5").Value = "8.815170815523743e-14"
$ws.Range("P5").Value = "2.660447832769735"
$ws.Range("Q5").Value = "[2.4717635893108882, 2.8491320762285817]"
$ws.Range("T5").Value = "56.54883145810828"
$ws.Range("U5").Value = "[51.05396325810985, 62.04369965810671]"
$ws.Range("X5").Value = "13.32468468468478"
$ws.Range("Y5").Value = "12.63069069069078"
$ws.Range("Z5").Value = "14.01867867867878"
$ws.Range("F6").Value = "23.11000000000017"
$ws.Range("H6").Value = "2.220446049250313e-16"
$ws.Range("I6").Value = "2.220446049250313e-16"
$ws.Range("L6").Value = "42.10924878139178"
$ws.Range("M6").Value = "[34.58771091219568, 49.630786650587886]"
$ws.Range("N6").Value = "1.06581410364015e-14"
$ws.Range("O6").Value = "1.06581410364015e-14"
$ws.Range("P6").Value = "2.949763672739966"
$ws.Range("Q6").Value = "[2.76107942928112, 3.1384479161988126]"
$ws.Range("T6").Value = "52.12845455118777"
$ws.Range("U6").Value = "[47.53621464252442, 56.720694459851124]"
$ws.Range("X6").Value = "12.26056056056065"
$ws.Range("Y6").Value = "11.56656656656665"
$ws.Range("Z6").Value = "12.95455455455465"
$ws.Range("F7").Value = "23.11000000000017"
$ws.Range("L7").Value = "51.08717450732755"
$ws.Range("M7").Value = "[43.89608973923394, 58.27825927542115]"
$ws.Range("P7").Value = "3.138447916198813"
$ws.Range("Q7").Value = "[2.974921571867812, 3.3019742605298132]"
$ws.Range("T7").Value = "57.34746283348922"
$ws.Range("U7").Value = "[52.3665685633469, 62.32835710363155]"
$ws.Range("X7").Value = "11.56656656656665"
$ws.Range("Y7").Value = "10.96510510510519"
$ws.Range("Z7").Value = "12.16802802802812"
